$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31 holds the F-series fuse designators (part 307010002).
# Two fuses (F74, F75) were removed from the design, per the final version.
$ws.Range("A31").Value = "F2, F5, F6, F9, F10, F15, F16, F19, F20, F26, F27, F28, F29, F40, F41, F42, F43, F46, F47, F48, F49, F50, F51, F54, F55, F56, F57, F59, F60, F62, F76, F78, F79, F82, F83"
$ws.Range("C31").Value = 35

$ws.Range("C32").Select() | Out-Null
